$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Anxa1"
$ws.Cells.Item(2, 3).Value = "Fpr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 22.61022533333333
$ws.Cells.Item(2, 8).Value = 67.830676
$ws.Cells.Item(2, 9).Value = 0.04352672200082041
$ws.Cells.Item(2, 10).Value = 0.04795217939334551
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.050657
$ws.Cells.Item(2, 14).Value = 0.151971
$ws.Cells.Item(2, 15).Value = 0.002123644810648064
$ws.Cells.Item(2, 16).Value = 0.002123644810648064
$ws.Cells.Item(2, 17).Value = 1.145366184710667
$ws.Cells.Item(2, 18).Value = 10.308295662396
$ws.Cells.Item(2, 19).Value = 0.00009243529730156317
$ws.Cells.Item(2, 20).Value = 0.0001018333969279432

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Anxa1"
$ws.Cells.Item(3, 3).Value = "Fpr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 22.61022533333333
$ws.Cells.Item(3, 8).Value = 67.830676
$ws.Cells.Item(3, 9).Value = 0.04352672200082041
$ws.Cells.Item(3, 10).Value = 0.04795217939334551
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.006255666666667
$ws.Cells.Item(3, 14).Value = 3.018767
$ws.Cells.Item(3, 15).Value = 0.0421842909114609
$ws.Cells.Item(3, 16).Value = 0.0421842909114609
$ws.Cells.Item(3, 17).Value = 22.75166736627689
$ws.Cells.Item(3, 18).Value = 204.765006296492
$ws.Cells.Item(3, 19).Value = 0.001836143903304893
$ws.Cells.Item(3, 20).Value = 0.002022828685367448

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Anxa1"
$ws.Cells.Item(4, 3).Value = "Fpr2"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 22.61022533333333
$ws.Cells.Item(4, 8).Value = 67.830676
$ws.Cells.Item(4, 9).Value = 0.04352672200082041
$ws.Cells.Item(4, 10).Value = 0.04795217939334551
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 18.813815
$ws.Cells.Item(4, 14).Value = 56.441445
$ws.Cells.Item(4, 15).Value = 0.7887135162611822
$ws.Cells.Item(4, 16).Value = 0.7887135162611822
$ws.Cells.Item(4, 17).Value = 425.3845965296467
$ws.Cells.Item(4, 18).Value = 3828.46136876682
$ws.Cells.Item(4, 19).Value = 0.03433011396059003
$ws.Cells.Item(4, 20).Value = 0.03782053202171254

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Anxa1"
$ws.Cells.Item(5, 3).Value = "Fpr2"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 22.61022533333333
$ws.Cells.Item(5, 8).Value = 67.830676
$ws.Cells.Item(5, 9).Value = 0.04352672200082041
$ws.Cells.Item(5, 10).Value = 0.04795217939334551
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 3.983073
$ws.Cells.Item(5, 14).Value = 11.949219
$ws.Cells.Item(5, 15).Value = 0.1669785480167087
$ws.Cells.Item(5, 16).Value = 0.1669785480167088
$ws.Cells.Item(5, 17).Value = 90.058178049116
$ws.Cells.Item(5, 18).Value = 810.523602442044
$ws.Cells.Item(5, 19).Value = 0.007268028839623924
$ws.Cells.Item(5, 20).Value = 0.008006985289337576

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Anxa1"
$ws.Cells.Item(6, 3).Value = "Fpr2"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 178.3379163333334
$ws.Cells.Item(6, 8).Value = 535.0137490000001
$ws.Cells.Item(6, 9).Value = 0.3433165654922813
$ws.Cells.Item(6, 10).Value = 0.3782223144872436
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.050657
$ws.Cells.Item(6, 14).Value = 0.151971
$ws.Cells.Item(6, 15).Value = 0.002123644810648064
$ws.Cells.Item(6, 16).Value = 0.002123644810648064
$ws.Cells.Item(6, 17).Value = 9.034063827697668
$ws.Cells.Item(6, 18).Value = 81.30657444927901
$ws.Cells.Item(6, 19).Value = 0.0007290824427171994
$ws.Cells.Item(6, 20).Value = 0.000803209855432135

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Anxa1"
$ws.Cells.Item(7, 3).Value = "Fpr2"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 178.3379163333334
$ws.Cells.Item(7, 8).Value = 535.0137490000001
$ws.Cells.Item(7, 9).Value = 0.3433165654922813
$ws.Cells.Item(7, 10).Value = 0.3782223144872436
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.006255666666667
$ws.Cells.Item(7, 14).Value = 3.018767
$ws.Cells.Item(7, 15).Value = 0.0421842909114609
$ws.Cells.Item(7, 16).Value = 0.0421842909114609
$ws.Cells.Item(7, 17).Value = 179.4535388919426
$ws.Cells.Item(7, 18).Value = 1615.081850027483
$ws.Cells.Item(7, 19).Value = 0.01448256587345001
$ws.Cells.Item(7, 20).Value = 0.01595504014353594

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Anxa1"
$ws.Cells.Item(8, 3).Value = "Fpr2"
$ws.Cells.Item(8, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 178.3379163333334
$ws.Cells.Item(8, 8).Value = 535.0137490000001
$ws.Cells.Item(8, 9).Value = 0.3433165654922813
$ws.Cells.Item(8, 10).Value = 0.3782223144872436
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 18.813815
$ws.Cells.Item(8, 14).Value = 56.441445
$ws.Cells.Item(8, 15).Value = 0.7887135162611822
$ws.Cells.Item(8, 16).Value = 0.7887135162611822
$ws.Cells.Item(8, 17).Value = 3355.216565380812
$ws.Cells.Item(8, 18).Value = 30196.94908842731
$ws.Cells.Item(8, 19).Value = 0.2707784155601297
$ws.Cells.Item(8, 20).Value = 0.2983090515876766

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Anxa1"
$ws.Cells.Item(9, 3).Value = "Fpr2"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 178.3379163333334
$ws.Cells.Item(9, 8).Value = 535.0137490000001
$ws.Cells.Item(9, 9).Value = 0.3433165654922813
$ws.Cells.Item(9, 10).Value = 0.3782223144872436
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.983073
$ws.Cells.Item(9, 14).Value = 11.949219
$ws.Cells.Item(9, 15).Value = 0.1669785480167087
$ws.Cells.Item(9, 16).Value = 0.1669785480167088
$ws.Cells.Item(9, 17).Value = 710.3329394235591
$ws.Cells.Item(9, 18).Value = 6392.996454812032
$ws.Cells.Item(9, 19).Value = 0.05732650161598443
$ws.Cells.Item(9, 20).Value = 0.06315501290059894

# Row 10
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Anxa1"
$ws.Cells.Item(10, 3).Value = "Fpr2"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 90.63663000000001
$ws.Cells.Item(10, 8).Value = 271.90989
$ws.Cells.Item(10, 9).Value = 0.174483683330882
$ws.Cells.Item(10, 10).Value = 0.1922238225092264
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.050657
$ws.Cells.Item(10, 14).Value = 0.151971
$ws.Cells.Item(10, 15).Value = 0.002123644810648064
$ws.Cells.Item(10, 16).Value = 0.002123644810648064
$ws.Cells.Item(10, 17).Value = 4.59137976591
$ws.Cells.Item(10, 18).Value = 41.32241789319
$ws.Cells.Item(10, 19).Value = 0.0003705413686483877
$ws.Cells.Item(10, 20).Value = 0.0004082151231546532

# Row 11
$ws.Cells.Item(11, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 2).Value = "Anxa1"
$ws.Cells.Item(11, 3).Value = "Fpr2"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 90.63663000000001
$ws.Cells.Item(11, 8).Value = 271.90989
$ws.Cells.Item(11, 9).Value = 0.174483683330882
$ws.Cells.Item(11, 10).Value = 0.1922238225092264
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.006255666666667
$ws.Cells.Item(11, 14).Value = 3.018767
$ws.Cells.Item(11, 15).Value = 0.0421842909114609
$ws.Cells.Item(11, 16).Value = 0.0421842909114609
$ws.Cells.Item(11, 17).Value = 91.20362254507
$ws.Cells.Item(11, 18).Value = 820.8326029056301
$ws.Cells.Item(11, 19).Value = 0.007360470456933148
$ws.Cells.Item(11, 20).Value = 0.008108825648842233

# Row 12
$ws.Cells.Item(12, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 2).Value = "Anxa1"
$ws.Cells.Item(12, 3).Value = "Fpr2"
$ws.Cells.Item(12, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 90.63663000000001
$ws.Cells.Item(12, 8).Value = 271.90989
$ws.Cells.Item(12, 9).Value = 0.174483683330882
$ws.Cells.Item(12, 10).Value = 0.1922238225092264
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 18.813815
$ws.Cells.Item(12, 14).Value = 56.441445
$ws.Cells.Item(12, 15).Value = 0.7887135162611822
$ws.Cells.Item(12, 16).Value = 0.7887135162611822
$ws.Cells.Item(12, 17).Value = 1705.22078904345
$ws.Cells.Item(12, 18).Value = 15346.98710139105
$ws.Cells.Item(12, 19).Value = 0.1376176394101026
$ws.Cells.Item(12, 20).Value = 0.1516095269604174

# Row 13
$ws.Cells.Item(13, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13, 2).Value = "Anxa1"
$ws.Cells.Item(13, 3).Value = "Fpr2"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 90.63663000000001
$ws.Cells.Item(13, 8).Value = 271.90989
$ws.Cells.Item(13, 9).Value = 0.174483683330882
$ws.Cells.Item(13, 10).Value = 0.1922238225092264
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 3.983073
$ws.Cells.Item(13, 14).Value = 11.949219
$ws.Cells.Item(13, 15).Value = 0.1669785480167087
$ws.Cells.Item(13, 16).Value = 0.1669785480167088
$ws.Cells.Item(13, 17).Value = 361.01231376399
$ws.Cells.Item(13, 18).Value = 3249.11082387591
$ws.Cells.Item(13, 19).Value = 0.02913503209519789
$ws.Cells.Item(13, 20).Value = 0.03209725477681217

# Row 14
$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "Anxa1"
$ws.Cells.Item(14, 3).Value = "Fpr2"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 143.820236
$ws.Cells.Item(14, 8).Value = 287.640472
$ws.Cells.Item(14, 9).Value = 0.2768669192002915
$ws.Cells.Item(14, 10).Value = 0.2033443911738485
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.050657
$ws.Cells.Item(14, 14).Value = 0.151971
$ws.Cells.Item(14, 15).Value = 0.002123644810648064
$ws.Cells.Item(14, 16).Value = 0.002123644810648064
$ws.Cells.Item(14, 17).Value = 7.285501695052
$ws.Cells.Item(14, 18).Value = 43.713010170312
$ws.Cells.Item(14, 19).Value = 0.0005879669961998159
$ws.Cells.Item(14, 20).Value = 0.0004318312610907333

# Row 15
$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "Anxa1"
$ws.Cells.Item(15, 3).Value = "Fpr2"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 143.820236
$ws.Cells.Item(15, 8).Value = 287.640472
$ws.Cells.Item(15, 9).Value = 0.2768669192002915
$ws.Cells.Item(15, 10).Value = 0.2033443911738485
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 1.006255666666667
$ws.Cells.Item(15, 14).Value = 3.018767
$ws.Cells.Item(15, 15).Value = 0.0421842909114609
$ws.Cells.Item(15, 16).Value = 0.0421842909114609
$ws.Cells.Item(15, 17).Value = 144.7199274563373
$ws.Cells.Item(15, 18).Value = 868.3195647380239
$ws.Cells.Item(15, 19).Value = 0.01167943466330503
$ws.Cells.Item(15, 20).Value = 0.008577938952491527

# Row 16
$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "Anxa1"
$ws.Cells.Item(16, 3).Value = "Fpr2"
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 143.820236
$ws.Cells.Item(16, 8).Value = 287.640472
$ws.Cells.Item(16, 9).Value = 0.2768669192002915
$ws.Cells.Item(16, 10).Value = 0.2033443911738485
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 18.813815
$ws.Cells.Item(16, 14).Value = 56.441445
$ws.Cells.Item(16, 15).Value = 0.7887135162611822
$ws.Cells.Item(16, 16).Value = 0.7887135162611822
$ws.Cells.Item(16, 17).Value = 2705.80731336034
$ws.Cells.Item(16, 18).Value = 16234.84388016204
$ws.Cells.Item(16, 19).Value = 0.2183686813788625
$ws.Cells.Item(16, 20).Value = 0.1603804697747153

# Row 17
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Anxa1"
$ws.Cells.Item(17, 3).Value = "Fpr2"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 143.820236
$ws.Cells.Item(17, 8).Value = 287.640472
$ws.Cells.Item(17, 9).Value = 0.2768669192002915
$ws.Cells.Item(17, 10).Value = 0.2033443911738485
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 3.983073
$ws.Cells.Item(17, 14).Value = 11.949219
$ws.Cells.Item(17, 15).Value = 0.1669785480167087
$ws.Cells.Item(17, 16).Value = 0.1669785480167088
$ws.Cells.Item(17, 17).Value = 572.8464988652279
$ws.Cells.Item(17, 18).Value = 3437.078993191368
$ws.Cells.Item(17, 19).Value = 0.0462308361619241
$ws.Cells.Item(17, 20).Value = 0.03395415118555087

# Row 18
$ws.Cells.Item(18, 1).Value = "Resolving-Mac"
$ws.Cells.Item(18, 2).Value = "Anxa1"
$ws.Cells.Item(18, 3).Value = "Fpr2"
$ws.Cells.Item(18, 4).Value = "ECs"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 84.051186
$ws.Cells.Item(18, 8).Value = 252.153558
$ws.Cells.Item(18, 9).Value = 0.1618061099757246
$ws.Cells.Item(18, 10).Value = 0.1782572924363359
$ws.Cells.Item(18, 11).Value = 1
$ws.Cells.Item(18, 12).Value = 0.3333333333333333
$ws.Cells.Item(18, 13).Value = 0.050657
$ws.Cells.Item(18, 14).Value = 0.151971
$ws.Cells.Item(18, 15).Value = 0.002123644810648064
$ws.Cells.Item(18, 16).Value = 0.002123644810648064
$ws.Cells.Item(18, 17).Value = 4.257780929202
$ws.Cells.Item(18, 18).Value = 38.320028362818
$ws.Cells.Item(18, 19).Value = 0.0003436187057810976
$ws.Cells.Item(18, 20).Value = 0.0003785551740425991

# Row 19
$ws.Cells.Item(19, 1).Value = "Resolving-Mac"
$ws.Cells.Item(19, 2).Value = "Anxa1"
$ws.Cells.Item(19, 3).Value = "Fpr2"
$ws.Cells.Item(19, 4).Value = "FAPs"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 84.051186
$ws.Cells.Item(19, 8).Value = 252.153558
$ws.Cells.Item(19, 9).Value = 0.1618061099757246
$ws.Cells.Item(19, 10).Value = 0.1782572924363359
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 1.006255666666667
$ws.Cells.Item(19, 14).Value = 3.018767
$ws.Cells.Item(19, 15).Value = 0.0421842909114609
$ws.Cells.Item(19, 16).Value = 0.0421842909114609
$ws.Cells.Item(19, 17).Value = 84.576982202554
$ws.Cells.Item(19, 18).Value = 761.192839822986
$ws.Cells.Item(19, 19).Value = 0.006825676014467804
$ws.Cells.Item(19, 20).Value = 0.007519657481223752

# Row 20
$ws.Cells.Item(20, 1).Value = "Resolving-Mac"
$ws.Cells.Item(20, 2).Value = "Anxa1"
$ws.Cells.Item(20, 3).Value = "Fpr2"
$ws.Cells.Item(20, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 84.051186
$ws.Cells.Item(20, 8).Value = 252.153558
$ws.Cells.Item(20, 9).Value = 0.1618061099757246
$ws.Cells.Item(20, 10).Value = 0.1782572924363359
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 18.813815
$ws.Cells.Item(20, 14).Value = 56.441445
$ws.Cells.Item(20, 15).Value = 0.7887135162611822
$ws.Cells.Item(20, 16).Value = 0.7887135162611822
$ws.Cells.Item(20, 17).Value = 1581.32346393459
$ws.Cells.Item(20, 18).Value = 14231.91117541131
$ws.Cells.Item(20, 19).Value = 0.1276186659514973
$ws.Cells.Item(20, 20).Value = 0.1405939359166603

# Row 21
$ws.Cells.Item(21, 1).Value = "Resolving-Mac"
$ws.Cells.Item(21, 2).Value = "Anxa1"
$ws.Cells.Item(21, 3).Value = "Fpr2"
$ws.Cells.Item(21, 4).Value = "Resolving-Mac"
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 84.051186
$ws.Cells.Item(21, 8).Value = 252.153558
$ws.Cells.Item(21, 9).Value = 0.1618061099757246
$ws.Cells.Item(21, 10).Value = 0.1782572924363359
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 3.983073
$ws.Cells.Item(21, 14).Value = 11.949219
$ws.Cells.Item(21, 15).Value = 0.1669785480167087
$ws.Cells.Item(21, 16).Value = 0.1669785480167088
$ws.Cells.Item(21, 17).Value = 334.782009574578
$ws.Cells.Item(21, 18).Value = 3013.038086171202
$ws.Cells.Item(21, 19).Value = 0.02701814930397839
$ws.Cells.Item(21, 20).Value = 0.02976514386440921

